$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "1.001", "30.707.40") that must
# remain plain text, matching the source data exactly. Temporarily force a text
# number format while assigning the value, then restore the default "Normal"
# style so no stray style index is left on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.707.40'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.891.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.35%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4889'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2942'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06689'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.889.85'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '89.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.023'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.649.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007911'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.06'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.133.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.749'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '192.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.085'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.325'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.36'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.833'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.402'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.277'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09039'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05213'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('E35').Value = '  -4.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.685'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01828'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.51%  '
$ws.Range('E38').Value = '  -2.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9249'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.051'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4416'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.752'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1344'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.374'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4160'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05833'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.716'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.413'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.30%  '
